# Update "Des Scheduled Flights vs actual.xlsx"
# Append 23 new daily rows (2021-10-03 .. 2021-10-25) below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data: Date, Scheduled flights (B), Tracked/actual flights (C)
$dates = @(
  "2021-10-03","2021-10-04","2021-10-05","2021-10-06","2021-10-07",
  "2021-10-08","2021-10-09","2021-10-10","2021-10-11","2021-10-12",
  "2021-10-13","2021-10-14","2021-10-15","2021-10-16","2021-10-17",
  "2021-10-18","2021-10-19","2021-10-20","2021-10-21","2021-10-22",
  "2021-10-23","2021-10-24","2021-10-25"
)
$scheduled = @(60,68,66,67,76,76,63,67,66,65,61,74,69,57,60,80,81,77,76,89,61,67,78)
$tracked   = @(59,68,65,66,71,70,60,61,61,61,58,70,68,55,57,77,75,70,68,75,59,62,70)

$firstNewRow = 546
$lastDataRow = 545
$lastNewRow = $firstNewRow + $dates.Length - 1

# Copy formatting (styles + number formats) from the last existing data row
# down across the whole new block, matching the prior row-by-row layout.
$ws.Range("A$($lastDataRow):D$($lastDataRow)").Copy() | Out-Null
$ws.Range("A$($firstNewRow):D$($lastNewRow)").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $dates.Length; $i++) {
  $r = $firstNewRow + $i
  $ws.Range("A$r").Value = $dates[$i]
  $ws.Range("B$r").Value = $scheduled[$i]
  $ws.Range("C$r").Value = $tracked[$i]
}

$ws.Range("D$($firstNewRow):D$($lastNewRow)").FormulaR1C1 = "=RC[-1]/RC[-2]"

# Match the new selection / active cell from the edit.
$ws.Range("I557").Select() | Out-Null
